$wb = $excel.ActiveWorkbook

# The two worksheets "展览" and "全部类型" contain identical data tables
# with a "想去人数" (interested count) column in column F for rows 3-5.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 20
    $ws.Range("F4").Value = 41
    $ws.Range("F5").Value = 7
}
